# Add a new "assignmentSearchBox" worksheet with data-table scenarios.
$wb = $excel.ActiveWorkbook

# Insert the new sheet after the last existing sheet, so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "assignmentSearchBox"

# Header row (B1:D1 - A1 "datakey" matches the other sheets' column A header)
$newSheet.Range("A1").Value = "datakey"

# Match the "datakey" header style used on the other sheets (bold Menlo font).
$sheet2 = $wb.Worksheets.Item(2)
$sheet2.Range("A1").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$newSheet.Range("B1").Value = "scenario"
$newSheet.Range("C1").Value = "feildName"
$newSheet.Range("D1").Value = "feildValue"

# Column A - scenario keys (data-table row labels)
$newSheet.Range("A2").Value = "valid_assgn_name"
$newSheet.Range("B2").Value = "valid"
$newSheet.Range("B3").Value = "invalid"
$newSheet.Range("A3").Value = "invalid_assgn_name"
$newSheet.Range("A4").Value = "valid_assgn_desc"
$newSheet.Range("A5").Value = "invalid_assgn_desc"
$newSheet.Range("A6").Value = "valid_assgn_duedate"
$newSheet.Range("A7").Value = "invalid_assgn_duedate"
$newSheet.Range("A8").Value = "valid_assgn_grade"
$newSheet.Range("A9").Value = "invalid_assgn_grade"

# Column B - valid/invalid scenario flag for remaining rows
$newSheet.Range("B4").Value = "valid"
$newSheet.Range("B5").Value = "invalid"
$newSheet.Range("B6").Value = "valid"
$newSheet.Range("B7").Value = "invalid"
$newSheet.Range("B8").Value = "valid"
$newSheet.Range("B9").Value = "invalid"

# Column C - field name being tested
$newSheet.Range("C2").Value = "name"
$newSheet.Range("C3").Value = "name"
$newSheet.Range("C4").Value = "desc"
$newSheet.Range("C5").Value = "desc"
$newSheet.Range("C6").Value = "dueDate"
$newSheet.Range("C7").Value = "dueDate"
$newSheet.Range("C8").Value = "grade"
$newSheet.Range("C9").Value = "grade"

# Column D - field value being tested
$newSheet.Range("D2").Value = "sdet-selenium"
$newSheet.Range("D3").Value = "xyz"
$newSheet.Range("D4").Value = "selenium"
$newSheet.Range("D5").Value = "xyz"
$newSheet.Range("D6").Value = 45224
$newSheet.Range("D6").NumberFormat = "mm-dd-yy"
$newSheet.Range("D7").Value = 9800
$newSheet.Range("D8").Value = "A"
$newSheet.Range("D9").Value = 123

$newSheet.Columns.Item(1).ColumnWidth = 15.666666666666666

# Make the new sheet the active / selected tab
$newSheet.Range("D1").Select()
$newSheet.Activate()
